# Generate Report for Handback
# Fills in "Latest Target File" (I) and "Latest Handback File" (J) /
# "Latest Handback DateTime" (K) columns on the zh-cn and de-de sheets,
# updates the Status column, and widens a few columns to fit the new
# (longer) content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f8b97ba350721a24922b850069a51c020cff0fe/e2e/62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$targetMdDisplay = "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f8b97ba350721a24922b850069a51c020cff0fe/e2e/62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$otherMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f8b97ba350721a24922b850069a51c020cff0fe/e2e/ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: Status, Latest Target File (hyperlink), Latest Handback
# File, Latest Handback DateTime
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = $statusHandedBack
$ws2.Range("C3").Value = $statusHandedBack

# Rebuild the hyperlinks collection in row-major order (A2, I2, A3, I3)
# so the new "Latest Target File" links land in between the pre-existing
# ones, matching how the handback report generator lays them out.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $sourceMdUrl, "", "", "62ad09a9-8825-480a-b9c6-c9b050f56804.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $otherMdUrl, "", "", "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null

$ws2.Range("J2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.zh-cn.xlf"
$ws2.Range("J3").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.zh-cn.xlf"

$ws2.Range("K2").Value = "2016-08-28 23:01:52"
$ws2.Range("K3").Value = "2016-08-28 23:01:52"

# ---------------------------------------------------------------------
# de-de sheet: Status, Latest Target File (hyperlink), Latest Handback
# File, Latest Handback DateTime
# ---------------------------------------------------------------------
$ws3.Range("C2").Value = $statusHandedBack
$ws3.Range("C3").Value = $statusHandedBack

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $sourceMdUrl, "", "", "62ad09a9-8825-480a-b9c6-c9b050f56804.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $otherMdUrl, "", "", "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null

$ws3.Range("J2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.de-de.xlf"
$ws3.Range("J3").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.de-de.xlf"

$ws3.Range("K2").Value = "2016-08-28 23:01:59"
$ws3.Range("K3").Value = "2016-08-28 23:01:59"

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn/de-de status columns pick up the new shared
# string automatically since they already reference the same text; make
# sure they're explicitly set too so the shared string gets reused
# instead of duplicated.
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $statusHandedBack
$ws1.Range("F2").Value = $statusHandedBack
$ws1.Range("E3").Value = $statusHandedBack
$ws1.Range("F3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Column widths: widen columns that now hold longer handback file names
# / URLs. ColumnWidth is rounded to whole pixels internally, so we pick
# the closest representable value to the target.
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664
